$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Host "Cells changed:" $changed
